$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logical Operators")

# Row 2: standalone formulas (not shared)
$ws.Range("D2").Formula = '=IF(C2>=60, "PASS", "FAIL")'
$ws.Range("E2").Formula = '=IF(C2>=90, "A", IF(C2>=80, "B", IF(C2>=70, "C", IF(C2 >=60, "D", "F"))))'
$ws.Range("F2").Formula = '=IF(OR(C2<60, C2>90), "Outlier", "Avg")'
$ws.Range("G2").Formula = '=IF(AND(B2="M",C2>95), "Male Achiever", IF(AND(B2="F",C2>95),"Female Achiever", "None"))'

# Rows 3:16: shared formulas across each column
$ws.Range("D3:D16").Formula = '=IF(C3>=60, "PASS", "FAIL")'
$ws.Range("E3:E16").Formula = '=IF(C3>=90, "A", IF(C3>=80, "B", IF(C3>=70, "C", IF(C3 >=60, "D", "F"))))'
$ws.Range("F3:F16").Formula = '=IF(OR(C3<60, C3>90), "Outlier", "Avg")'
$ws.Range("G3:G16").Formula = '=IF(AND(B3="M",C3>95), "Male Achiever", IF(AND(B3="F",C3>95),"Female Achiever", "None"))'

# Update the selected cell in the sheet view
$ws.Range("F9").Select()
